# Rewrite the META sheet so that the consumption data (p5/p6 producers and the
# residential_-5 / commercial_-8 consumers) is generated straight from this
# meta sheet instead of relying on the separate CONSUMPTION file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- wipe out the old data area (everything below the header row) ---------
$ws.Range("A2:H59").Clear()

# --- coordinates block (rows 2-9) ------------------------------------------
$coords = @(
    @("p1", "coordinates", "latitude,longitude", 43.6,               -116.2),
    @("p2", "coordinates", "latitude,longitude", 43.5,               -115.3),
    @("p3", "coordinates", "latitude,longitude", 51.6189802813036,   5.71901122838461),
    @("p4", "coordinates", "latitude,longitude", 51.5,               4.2300000000000004),
    @("p5", "coordinates", "latitude,longitude", 51.8902154304227,   4.1223949882912398),
    @("p6", "coordinates", "latitude,longitude", 51.8902154304227,   -40),
    @("c1", "coordinates", "latitude,longitude", 52.343337996743898, 4.9358207840021304),
    @("c2", "coordinates", "latitude,longitude", 53.198799420490097, 6.5688442594149397)
)
$r = 2
foreach ($row in $coords) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 4).Value2 = $row[2]
    $ws.Cells.Item($r, 5).Value2 = $row[3]
    $ws.Cells.Item($r, 6).Value2 = $row[4]
    $r = $r + 1
}

# --- type block (rows 10-17) ------------------------------------------------
$types = @(
    @("p1", "solar"),
    @("p2", "wind"),
    @("p3", "solar"),
    @("p4", "wind"),
    @("p5", "biomass"),
    @("p6", "solar"),
    @("c1", "residential_-5"),
    @("c2", "commercial_-8")
)
foreach ($row in $types) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = "type"
    $ws.Cells.Item($r, 5).Value2 = $row[1]
    $r = $r + 1
}

# --- capacity block (rows 18-25) --------------------------------------------
$capacities = @(
    @("p1", 30),
    @("p2", 30),
    @("p3", 30),
    @("p4", 30),
    @("p5", 30),
    @("p6", 30),
    @("c1", 100),
    @("c2", 50)
)
foreach ($row in $capacities) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = "capacity"
    $ws.Cells.Item($r, 3).Value2 = "GWh"
    $ws.Cells.Item($r, 4).Value2 = "year average"
    $valueCell = $ws.Cells.Item($r, 5)
    $valueCell.Value2 = $row[1]
    if ($row[0] -ne "c1" -and $row[0] -ne "c2") {
        $valueCell.NumberFormat = "#,##0"
    }
    $r = $r + 1
}

# --- DC line / transmission block (rows 26-28) ------------------------------
$ws.Cells.Item($r, 1).Value2 = "t1"
$ws.Cells.Item($r, 2).Value2 = "efficiency"
$ws.Cells.Item($r, 3).Value2 = "%/1000 km"
$ws.Cells.Item($r, 5).Value2 = 0.997
$r = $r + 1

$ws.Cells.Item($r, 1).Value2 = "t1"
$ws.Cells.Item($r, 2).Value2 = "type"
$ws.Cells.Item($r, 5).Value2 = "DC line"
$r = $r + 1

$ws.Cells.Item($r, 1).Value2 = "t1"
$ws.Cells.Item($r, 2).Value2 = "conversion"
$ws.Cells.Item($r, 4).Value2 = "conversion losses"
$ws.Cells.Item($r, 5).Value2 = 0.98499999999999999
$r = $r + 1

# --- solar constants block (rows 29-32) -------------------------------------
$constants = @(
    @("d0", 23.45),
    @("a0", 0.42370000000000002),
    @("a1", 0.50549999999999995),
    @("k",  0.27110000000000001)
)
foreach ($row in $constants) {
    $ws.Cells.Item($r, 1).Value2 = "constant"
    $ws.Cells.Item($r, 2).Value2 = "solar"
    $ws.Cells.Item($r, 4).Value2 = $row[0]
    $ws.Cells.Item($r, 5).Value2 = $row[1]
    $r = $r + 1
}

# --- trailing formatted-but-empty cells in column E -------------------------
$ws.Range("E45:E48").NumberFormat = "#,##0"
$ws.Range("E62:E65").NumberFormat = "#,##0"

# --- selection / view state, matching the authored workbook ----------------
$ws.Range("E19:E23").Select()
